$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as literal text, preserving the default (unstyled)
# cell style even when the string looks like a number (e.g. "211.92", "1.00").
# Forcing NumberFormat="@" while assigning keeps Excel from coercing the text to
# a double, then Style="Normal" drops the temporary text format back off the cell
# so no stray "s" attribute is left on it.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value2 = $value
    $range.Style = "Normal"
}

$ws.Range('D2').Value2 = '27.731.07'
$ws.Range('E2').Value2 = '  -0.17%  '
$ws.Range('D3').Value2 = '1.635.79'
$ws.Range('E3').Value2 = '  +0.22%  '
Set-TextValue $ws.Range('D4') '0.999'
$ws.Range('E4').Value2 = '  -0.49%  '
Set-TextValue $ws.Range('D5') '211.92'
$ws.Range('E5').Value2 = '  -0.31%  '
Set-TextValue $ws.Range('D6') '0.523'
$ws.Range('E6').Value2 = '  -0.31%  '
Set-TextValue $ws.Range('D7') '0.999'
$ws.Range('E7').Value2 = '  -0.52%  '
Set-TextValue $ws.Range('D8') '23.25'
$ws.Range('E8').Value2 = '  +1.11%  '
Set-TextValue $ws.Range('D9') '0.264'
$ws.Range('E9').Value2 = '  +0.67%  '
Set-TextValue $ws.Range('D10') '0.0612'
$ws.Range('E10').Value2 = '  +0.07%  '
Set-TextValue $ws.Range('D11') '0.0863'
$ws.Range('E11').Value2 = '  -3.05%  '
$ws.Range('D12').Value2 = '1.865.27'
$ws.Range('E12').Value2 = '  +0.05%  '
$ws.Range('D13').Value2 = '1.635.48'
$ws.Range('E13').Value2 = '  +0.21%  '
Set-TextValue $ws.Range('D14') '4.05'
$ws.Range('E14').Value2 = '  +0.06%  '
Set-TextValue $ws.Range('D15') '0.562'
$ws.Range('E15').Value2 = '  +1.15%  '
Set-TextValue $ws.Range('D16') '65.28'
$ws.Range('E16').Value2 = '  +1.30%  '
$ws.Range('D17').Value2 = '27.661.16'
$ws.Range('E17').Value2 = '  -0.42%  '
Set-TextValue $ws.Range('D18') '230.07'
$ws.Range('E18').Value2 = '  -0.44%  '
$ws.Range('D19').Value2 = '0.0₃0720'
$ws.Range('E19').Value2 = '  -0.28%  '
Set-TextValue $ws.Range('D20') '7.62'
$ws.Range('E20').Value2 = '  +0.26%  '
Set-TextValue $ws.Range('D21') '1.00'
$ws.Range('E21').Value2 = '  -0.27%  '
Set-TextValue $ws.Range('D22') '10.70'
$ws.Range('E22').Value2 = '  +7.58%  '
$ws.Range('E23').Value2 = '  +2.03%  '
Set-TextValue $ws.Range('D24') '2.15'
$ws.Range('E24').Value2 = '  +2.96%  '
Set-TextValue $ws.Range('D25') '149.87'
$ws.Range('E25').Value2 = '  -0.03%  '
Set-TextValue $ws.Range('D26') '6.89'
$ws.Range('E26').Value2 = '  -0.49%  '
$ws.Range('B27').Value2 = 'EthereumClassic'
$ws.Range('C27').Value2 = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range('D27') '15.67'
$ws.Range('E27').Value2 = '  +0.44%  '
$ws.Range('B28').Value2 = 'Stellar'
$ws.Range('C28').Value2 = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range('D28') '0.111'
$ws.Range('E28').Value2 = '  -0.37%  '
Set-TextValue $ws.Range('D29') '0.999'
$ws.Range('E29').Value2 = '  -0.41%  '
$ws.Range('E30').Value2 = '  -0.22%  '
Set-TextValue $ws.Range('D31') '0.0482'
$ws.Range('E31').Value2 = '  -0.21%  '
Set-TextValue $ws.Range('D32') '3.29'
$ws.Range('E32').Value2 = '  -0.27%  '
$ws.Range('D33').Value2 = '1.467.60'
$ws.Range('E33').Value2 = '  +0.00%  '
$ws.Range('E34').Value2 = '  +0.10%  '
$ws.Range('E35').Value2 = '  +0.29%  '
$ws.Range('E36').Value2 = '  -1.90%  '
$ws.Range('B37').Value2 = 'TrustWalletToken'
$ws.Range('C37').Value2 = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range('D37') '0.930'
$ws.Range('E37').Value2 = '  +1.79%  '
$ws.Range('B38').Value2 = 'ARBITRUM'
$ws.Range('C38').Value2 = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range('D38') '0.880'
$ws.Range('E38').Value2 = '  +1.00%  '
$ws.Range('E39').Value2 = '  -0.15%  '
Set-TextValue $ws.Range('D40') '0.556'
$ws.Range('E40').Value2 = '  -1.53%  '
Set-TextValue $ws.Range('D41') '68.97'
$ws.Range('E41').Value2 = '  -0.50%  '
$ws.Range('B42').Value2 = 'WEMIXToken'
$ws.Range('C42').Value2 = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range('D42') '1.02'
$ws.Range('E42').Value2 = '  -0.52%  '
$ws.Range('B43').Value2 = 'PaxDollar'
$ws.Range('C43').Value2 = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue $ws.Range('D43') '1.00'
$ws.Range('E43').Value2 = '  -0.33%  '
$ws.Range('E44').Value2 = '  -0.02%  '
$ws.Range('E45').Value2 = '  -0.72%  '
$ws.Range('E46').Value2 = '  -0.76%  '
$ws.Range('D47').Value2 = '1.774.98'
$ws.Range('E47').Value2 = '  -0.05%  '
Set-TextValue $ws.Range('D48') '1.75'
$ws.Range('E48').Value2 = '  +3.35%  '
Set-TextValue $ws.Range('D49') '87.74'
$ws.Range('E49').Value2 = '  +2.17%  '
$ws.Range('E50').Value2 = '  +4.53%  '
Set-TextValue $ws.Range('D51') '0.0999'
$ws.Range('E51').Value2 = '  +0.90%  '
